$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell value while forcing text storage (avoids Excel
# auto-converting numeric-looking strings like "1.01" into real numbers),
# and resetting the style afterwards so no stray number-format style is left on the cell.
function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "25.827.49"
Set-TextValue $ws.Range("E2") "  +0.00%  "

Set-TextValue $ws.Range("D3") "1.637.10"
Set-TextValue $ws.Range("E3") "  +0.06%  "

Set-TextValue $ws.Range("D4") "1.01"
Set-TextValue $ws.Range("E4") "  +0.53%  "

Set-TextValue $ws.Range("D5") "216.18"
Set-TextValue $ws.Range("E5") "  +0.43%  "

Set-TextValue $ws.Range("E6") "  -0.08%  "

Set-TextValue $ws.Range("D7") "1.01"
Set-TextValue $ws.Range("E7") "  +0.45%  "

Set-TextValue $ws.Range("D8") "0.258"
Set-TextValue $ws.Range("E8") "  +0.46%  "

Set-TextValue $ws.Range("D9") "0.0640"
Set-TextValue $ws.Range("E9") "  -0.28%  "

Set-TextValue $ws.Range("D10") "19.67"
Set-TextValue $ws.Range("E10") "  -2.24%  "

Set-TextValue $ws.Range("D11") "0.0785"
Set-TextValue $ws.Range("E11") "  +0.74%  "

Set-TextValue $ws.Range("D12") "1.673.25"
Set-TextValue $ws.Range("E12") "  +2.18%  "

Set-TextValue $ws.Range("D13") "4.26"
Set-TextValue $ws.Range("E13") "  -0.65%  "

Set-TextValue $ws.Range("D14") "1.857.76"
Set-TextValue $ws.Range("E14") "  -0.24%  "

Set-TextValue $ws.Range("D15") "0.556"
Set-TextValue $ws.Range("E15") "  -0.81%  "

Set-TextValue $ws.Range("D16") "0.0₃0775"
Set-TextValue $ws.Range("E16") "  +1.47%  "

Set-TextValue $ws.Range("D17") "63.35"
Set-TextValue $ws.Range("E17") "  +0.48%  "

Set-TextValue $ws.Range("D18") "25.861.98"
Set-TextValue $ws.Range("E18") "  +0.01%  "

Set-TextValue $ws.Range("D19") "1.01"
Set-TextValue $ws.Range("E19") "  +0.32%  "

Set-TextValue $ws.Range("E20") "  +2.42%  "

Set-TextValue $ws.Range("D21") "194.29"
Set-TextValue $ws.Range("E21") "  -0.11%  "

Set-TextValue $ws.Range("D22") "9.97"
Set-TextValue $ws.Range("E22") "  +0.93%  "

Set-TextValue $ws.Range("D23") "6.19"
Set-TextValue $ws.Range("E23") "  +1.47%  "

Set-TextValue $ws.Range("E24") "  +0.45%  "

Set-TextValue $ws.Range("D25") "1.77"
Set-TextValue $ws.Range("E25") "  -0.41%  "

Set-TextValue $ws.Range("D26") "140.03"
Set-TextValue $ws.Range("E26") "  -0.29%  "

Set-TextValue $ws.Range("E27") "  -3.70%  "

Set-TextValue $ws.Range("D28") "6.86"
Set-TextValue $ws.Range("E28") "  +0.69%  "

Set-TextValue $ws.Range("D29") "15.63"
Set-TextValue $ws.Range("E29") "  +1.02%  "

Set-TextValue $ws.Range("D30") "1.25"
Set-TextValue $ws.Range("E30") "  +0.43%  "

Set-TextValue $ws.Range("E31") "  -0.15%  "

Set-TextValue $ws.Range("D32") "3.36"
Set-TextValue $ws.Range("E32") "  +1.74%  "

Set-TextValue $ws.Range("D33") "3.28"
Set-TextValue $ws.Range("E33") "  +1.46%  "

Set-TextValue $ws.Range("E34") "  +1.56%  "

Set-TextValue $ws.Range("D35") "2.40"
Set-TextValue $ws.Range("E35") "  +0.85%  "

Set-TextValue $ws.Range("D36") "0.898"
Set-TextValue $ws.Range("E36") "  -0.58%  "

Set-TextValue $ws.Range("E37") "  +0.22%  "

Set-TextValue $ws.Range("D38") "0.553"
Set-TextValue $ws.Range("E38") "  +0.00%  "

Set-TextValue $ws.Range("D39") "1.109.38"
Set-TextValue $ws.Range("E39") "  -1.45%  "

Set-TextValue $ws.Range("E40") "  +0.67%  "

Set-TextValue $ws.Range("D41") "1.01"
Set-TextValue $ws.Range("E41") "  +0.37%  "

Set-TextValue $ws.Range("D42") "5.59"
Set-TextValue $ws.Range("E42") "  +1.34%  "

Set-TextValue $ws.Range("D43") "0.806"
Set-TextValue $ws.Range("E43") "  +0.87%  "

Set-TextValue $ws.Range("D44") "99.61"
Set-TextValue $ws.Range("E44") "  +1.53%  "

Set-TextValue $ws.Range("D45") "0.0₆0110"
Set-TextValue $ws.Range("E45") "  -2.22%  "

Set-TextValue $ws.Range("D46") "55.23"
Set-TextValue $ws.Range("E46") "  -0.34%  "

Set-TextValue $ws.Range("D47") "2.43"
Set-TextValue $ws.Range("E47") "  +10.51%  "

Set-TextValue $ws.Range("B48") "EnergySwap"
Set-TextValue $ws.Range("C48") "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue $ws.Range("D48") "7.75"
Set-TextValue $ws.Range("E48") "  +0.12%  "

Set-TextValue $ws.Range("B49") "Mantle"
Set-TextValue $ws.Range("C49") "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
Set-TextValue $ws.Range("D49") "0.419"
Set-TextValue $ws.Range("E49") "  -1.71%  "

Set-TextValue $ws.Range("E50") "  +0.20%  "

Set-TextValue $ws.Range("E51") "  +0.35%  "

